# Auto-generated edit script updating the crypto price/volume/hour table
# Values are entered with a leading apostrophe so Excel stores them as
# literal text (matching the source data's inline-string cells) instead
# of auto-converting numeric-looking / percent-looking text into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''308.44'
$ws.Range("E2").Value = '''-0.87%'
$ws.Range("G2").Value = '''10'
$ws.Range("D3").Value = '''37.89'
$ws.Range("E3").Value = '''-3.54%'
$ws.Range("G3").Value = '''10'
$ws.Range("D4").Value = '''5.062'
$ws.Range("E4").Value = '''-2.10%'
$ws.Range("G4").Value = '''10'
$ws.Range("D5").Value = '''0.07891'
$ws.Range("E5").Value = '''-3.18%'
$ws.Range("G5").Value = '''10'
$ws.Range("D6").Value = '''2.027'
$ws.Range("E6").Value = '''0.86%'
$ws.Range("G6").Value = '''10'
$ws.Range("D7").Value = '''4.394'
$ws.Range("E7").Value = '''3.60%'
$ws.Range("G7").Value = '''10'
$ws.Range("D8").Value = '''8.270'
$ws.Range("E8").Value = '''1.50%'
$ws.Range("G8").Value = '''10'
$ws.Range("D9").Value = '''3.114'
$ws.Range("E9").Value = '''-7.78%'
$ws.Range("G9").Value = '''10'
$ws.Range("D10").Value = '''0.9320'
$ws.Range("E10").Value = '''0.59%'
$ws.Range("G10").Value = '''10'
$ws.Range("E11").Value = '''-7.48%'
$ws.Range("G11").Value = '''10'
$ws.Range("E12").Value = '''-2.73%'
$ws.Range("G12").Value = '''10'
$ws.Range("D13").Value = '''0.08731'
$ws.Range("E13").Value = '''-3.57%'
$ws.Range("G13").Value = '''10'
$ws.Range("D14").Value = '''0.03457'
$ws.Range("E14").Value = '''-1.67%'
$ws.Range("G14").Value = '''10'
$ws.Range("D15").Value = '''0.09655'
$ws.Range("E15").Value = '''-1.67%'
$ws.Range("G15").Value = '''10'
$ws.Range("D16").Value = '''0.001402'
$ws.Range("E16").Value = '''1.00%'
$ws.Range("G16").Value = '''10'
$ws.Range("D17").Value = '''0.006223'
$ws.Range("E17").Value = '''2.01%'
$ws.Range("G17").Value = '''10'
$ws.Range("D18").Value = '''3.593'
$ws.Range("E18").Value = '''-2.24%'
$ws.Range("G18").Value = '''10'
$ws.Range("D19").Value = '''0.3435'
$ws.Range("E19").Value = '''-0.62%'
$ws.Range("G19").Value = '''10'
$ws.Range("D20").Value = '''0.1293'
$ws.Range("E20").Value = '''-1.37%'
$ws.Range("G20").Value = '''10'
$ws.Range("D21").Value = '''5.032'
$ws.Range("E21").Value = '''8.37%'
$ws.Range("G21").Value = '''10'
$ws.Range("D22").Value = '''0.2527'
$ws.Range("G22").Value = '''10'
$ws.Range("D23").Value = '''0.04366'
$ws.Range("E23").Value = '''-0.26%'
$ws.Range("G23").Value = '''10'
$ws.Range("D24").Value = '''0.001237'
$ws.Range("E24").Value = '''0.55%'
$ws.Range("G24").Value = '''10'
$ws.Range("D25").Value = '''0.004633'
$ws.Range("E25").Value = '''-4.96%'
$ws.Range("G25").Value = '''10'
$ws.Range("E26").Value = '''176.48%'
$ws.Range("G26").Value = '''10'
$ws.Range("G27").Value = '''10'
$ws.Range("G28").Value = '''10'
$ws.Range("G29").Value = '''10'
$ws.Range("G30").Value = '''10'
$ws.Range("G31").Value = '''10'
$ws.Range("G32").Value = '''10'
$ws.Range("G33").Value = '''10'
$ws.Range("G34").Value = '''10'
$ws.Range("G35").Value = '''10'
$ws.Range("G36").Value = '''10'
$ws.Range("G37").Value = '''10'
$ws.Range("G38").Value = '''10'
$ws.Range("D39").Value = '''0.02210'
$ws.Range("E39").Value = '''3.21%'
$ws.Range("G39").Value = '''10'
$ws.Range("D40").Value = '''0.05044'
$ws.Range("E40").Value = '''-2.99%'
$ws.Range("G40").Value = '''10'
$ws.Range("D41").Value = '''0.007533'
$ws.Range("E41").Value = '''1.37%'
$ws.Range("G41").Value = '''10'
$ws.Range("D42").Value = '''0.01001'
$ws.Range("E42").Value = '''1.78%'
$ws.Range("G42").Value = '''10'
$ws.Range("D43").Value = '''0.1365'
$ws.Range("E43").Value = '''-0.19%'
$ws.Range("G43").Value = '''10'
$ws.Range("D44").Value = '''0.002041'
$ws.Range("E44").Value = '''-4.12%'
$ws.Range("G44").Value = '''10'
$ws.Range("D45").Value = '''0.008841'
$ws.Range("E45").Value = '''-10.46%'
$ws.Range("G45").Value = '''10'
$ws.Range("D46").Value = '''0.00006670'
$ws.Range("E46").Value = '''4.16%'
$ws.Range("G46").Value = '''10'
$ws.Range("D47").Value = '''0.00000000755'
$ws.Range("E47").Value = '''0.71%'
$ws.Range("G47").Value = '''10'
$ws.Range("D48").Value = '''0.003019'
$ws.Range("E48").Value = '''9.65%'
$ws.Range("G48").Value = '''10'
$ws.Range("D49").Value = '''0.001207'
$ws.Range("E49").Value = '''20.81%'
$ws.Range("G49").Value = '''10'
$ws.Range("D50").Value = '''0.00002114'
$ws.Range("E50").Value = '''0.71%'
$ws.Range("G50").Value = '''10'
$ws.Range("D51").Value = '''0.0002013'
$ws.Range("E51").Value = '''0.71%'
$ws.Range("G51").Value = '''10'
